$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Samuel" (row 6) to "Tote"
$ws.Range("A6").Value = "Tote"

# Delete the entire row 7 (previously "Francisco" row), shifting cells up
$ws.Rows("7:7").Delete()

# Update the active selection to match the target state
$ws.Range("B7").Select()
